$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells whose new values would otherwise be auto-parsed as numbers,
# so they are written back as text (matching the source data which stores them as strings).
$textCoercedCells = @(
    "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D15", "D19", "D21", "D23", "D25", "D26", "D27", "D28", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D42", "D45", "D47", "D48"
)
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '36.599.83'
$ws.Range("E2").Value = '  +3.17%  '
$ws.Range("D3").Value = '2.074.84'
$ws.Range("E3").Value = '  +10.01%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '247.44'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = '0.665'
$ws.Range("E6").Value = '  -3.10%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '45.19'
$ws.Range("E8").Value = '  +5.56%  '
$ws.Range("D9").Value = '60.70'
$ws.Range("E9").Value = '  +7.35%  '
$ws.Range("D10").Value = '0.366'
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("D11").Value = '0.0720'
$ws.Range("E11").Value = '  -4.12%  '
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '14.54'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").Value = '2.361.06'
$ws.Range("E14").Value = '  +9.07%  '
$ws.Range("D15").Value = '0.818'
$ws.Range("E15").Value = '  +3.31%  '
$ws.Range("D16").Value = '2.050.55'
$ws.Range("E16").Value = '  +8.42%  '
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").Value = '36.601.04'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").Value = '71.46'
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("D20").Value = '0.0₃0812'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("D21").Value = '237.53'
$ws.Range("E21").Value = '  -3.47%  '
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("D23").Value = '4.92'
$ws.Range("E23").Value = '  -4.76%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  -7.26%  '
$ws.Range("D26").Value = '169.05'
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("D27").Value = '20.17'
$ws.Range("E27").Value = '  +10.06%  '
$ws.Range("D28").Value = '8.80'
$ws.Range("E28").Value = '  +2.02%  '
$ws.Range("E29").Value = '  -8.65%  '
$ws.Range("E30").Value = '  -4.67%  '
$ws.Range("D31").Value = '21.69'
$ws.Range("E31").Value = '  +49.75%  '
$ws.Range("D32").Value = '4.36'
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").Value = '0.0582'
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("D34").Value = '0.0903'
$ws.Range("E34").Value = '  +22.02%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +1.19%  '
$ws.Range("D37").Value = '2.26'
$ws.Range("E37").Value = '  +16.99%  '
$ws.Range("D38").Value = '3.98'
$ws.Range("E38").Value = '  -6.73%  '
$ws.Range("D39").Value = '0.873'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").Value = '1.33'
$ws.Range("E40").Value = '  -10.01%  '
$ws.Range("E41").Value = '  +5.01%  '
$ws.Range("D42").Value = '96.88'
$ws.Range("E42").Value = '  -1.91%  '
$ws.Range("E43").Value = '  -6.39%  '
$ws.Range("E44").Value = '  +15.94%  '
$ws.Range("D45").Value = '16.15'
$ws.Range("E45").Value = '  -4.69%  '
$ws.Range("D46").Value = '1.319.67'
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").Value = '0.0819'
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("D48").Value = '2.82'
$ws.Range("E48").Value = '  +3.35%  '
$ws.Range("D49").Value = '2.246.76'
$ws.Range("E49").Value = '  +8.82%  '
$ws.Range("E50").Value = '  -5.59%  '
$ws.Range("E51").Value = '  +15.59%  '

# Restore the General (default) format on the cells we temporarily marked as text,
# so no stray formatting is introduced.
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).ClearFormats()
}
